# Adds missing headers (external_id, utm_id, utm_source, utm_medium,
# utm_campaign, utm_term, utm_content) to the pw_purchases sheet's headers
# table (rows 9-15), matching the "purchases_table" schema on sheet 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pw_purchases")

$rows = @(
    @{ Row = 9;  Name = "external_id";  Example = "#ioy2fcf" },
    @{ Row = 10; Name = "utm_id";       Example = "wvninnewn" },
    @{ Row = 11; Name = "utm_source";   Example = "youtube" },
    @{ Row = 12; Name = "utm_medium";   Example = "social" },
    @{ Row = 13; Name = "utm_campaign"; Example = "fireship" },
    @{ Row = 14; Name = "utm_term";     Example = "apple" },
    @{ Row = 15; Name = "utm_content";  Example = "logolink" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Name       # A: headers
    $ws.Cells.Item($row, 2).Value = "Yes"          # B: null_allowed
    $ws.Cells.Item($row, 3).Value = "string"       # C: type
    $ws.Cells.Item($row, 11).Value = $r.Example    # K: examples
}
